$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header row (row 1) cell values
$ws.Range("A1").Value = "question"
$ws.Range("B1").Value = "Option A"
$ws.Range("C1").Value = "Option B"
$ws.Range("D1").Value = "Option C"
$ws.Range("E1").Value = "Option D"

# Update the selected cell in the sheet view
$ws.Range("G2").Select()
